$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.Value = "'" + '28.420.50'
$c.Style = 'Normal'
$c = $ws.Range('E2')
$c.Value = "'" + '  +0.14%  '
$c.Style = 'Normal'
$c = $ws.Range('D3')
$c.Value = "'" + '1.569.90'
$c.Style = 'Normal'
$c = $ws.Range('E3')
$c.Value = "'" + '  -1.66%  '
$c.Style = 'Normal'
$c = $ws.Range('D5')
$c.Value = "'" + '211.74'
$c.Style = 'Normal'
$c = $ws.Range('E5')
$c.Value = "'" + '  -1.44%  '
$c.Style = 'Normal'
$c = $ws.Range('E6')
$c.Value = "'" + '  -1.00%  '
$c.Style = 'Normal'
$c = $ws.Range('E7')
$c.Value = "'" + '  +0.00%  '
$c.Style = 'Normal'
$c = $ws.Range('D8')
$c.Value = "'" + '45.47'
$c.Style = 'Normal'
$c = $ws.Range('E8')
$c.Value = "'" + '  +3.05%  '
$c.Style = 'Normal'
$c = $ws.Range('D9')
$c.Value = "'" + '24.02'
$c.Style = 'Normal'
$c = $ws.Range('E9')
$c.Value = "'" + '  -0.71%  '
$c.Style = 'Normal'
$c = $ws.Range('E10')
$c.Value = "'" + '  -2.15%  '
$c.Style = 'Normal'
$c = $ws.Range('E11')
$c.Value = "'" + '  -1.85%  '
$c.Style = 'Normal'
$c = $ws.Range('D12')
$c.Value = "'" + '0.0889'
$c.Style = 'Normal'
$c = $ws.Range('E12')
$c.Value = "'" + '  +0.14%  '
$c.Style = 'Normal'
$c = $ws.Range('D13')
$c.Value = "'" + '1.792.89'
$c.Style = 'Normal'
$c = $ws.Range('E13')
$c.Value = "'" + '  -1.76%  '
$c.Style = 'Normal'
$c = $ws.Range('D14')
$c.Value = "'" + '1.562.30'
$c.Style = 'Normal'
$c = $ws.Range('E14')
$c.Value = "'" + '  -2.09%  '
$c.Style = 'Normal'
$c = $ws.Range('E15')
$c.Value = "'" + '  -2.34%  '
$c.Style = 'Normal'
$c = $ws.Range('D16')
$c.Value = "'" + '28.420.62'
$c.Style = 'Normal'
$c = $ws.Range('E16')
$c.Value = "'" + '  +0.08%  '
$c.Style = 'Normal'
$c = $ws.Range('D17')
$c.Value = "'" + '3.66'
$c.Style = 'Normal'
$c = $ws.Range('E17')
$c.Value = "'" + '  -2.52%  '
$c.Style = 'Normal'
$c = $ws.Range('D18')
$c.Value = "'" + '62.15'
$c.Style = 'Normal'
$c = $ws.Range('E18')
$c.Value = "'" + '  -1.69%  '
$c.Style = 'Normal'
$c = $ws.Range('D19')
$c.Value = "'" + '227.74'
$c.Style = 'Normal'
$c = $ws.Range('E19')
$c.Value = "'" + '  -0.01%  '
$c.Style = 'Normal'
$c = $ws.Range('E20')
$c.Value = "'" + '  -2.64%  '
$c.Style = 'Normal'
$c = $ws.Range('E21')
$c.Value = "'" + '  -3.32%  '
$c.Style = 'Normal'
$c = $ws.Range('E22')
$c.Value = "'" + '  +0.03%  '
$c.Style = 'Normal'
$c = $ws.Range('E23')
$c.Value = "'" + '  -6.00%  '
$c.Style = 'Normal'
$c = $ws.Range('E25')
$c.Value = "'" + '  +7.11%  '
$c.Style = 'Normal'
$c = $ws.Range('D26')
$c.Value = "'" + '150.72'
$c.Style = 'Normal'
$c = $ws.Range('E26')
$c.Value = "'" + '  -0.93%  '
$c.Style = 'Normal'
$c = $ws.Range('D27')
$c.Value = "'" + '14.94'
$c.Style = 'Normal'
$c = $ws.Range('E27')
$c.Value = "'" + '  -1.84%  '
$c.Style = 'Normal'
$c = $ws.Range('D28')
$c.Value = "'" + '6.43'
$c.Style = 'Normal'
$c = $ws.Range('E28')
$c.Value = "'" + '  -2.58%  '
$c.Style = 'Normal'
$c = $ws.Range('E29')
$c.Value = "'" + '  -3.57%  '
$c.Style = 'Normal'
$c = $ws.Range('E30')
$c.Value = "'" + '  -0.03%  '
$c.Style = 'Normal'
$c = $ws.Range('D31')
$c.Value = "'" + '0.0489'
$c.Style = 'Normal'
$c = $ws.Range('E31')
$c.Value = "'" + '  +2.79%  '
$c.Style = 'Normal'
$c = $ws.Range('E32')
$c.Value = "'" + '  -4.19%  '
$c.Style = 'Normal'
$c = $ws.Range('E33')
$c.Value = "'" + '  -1.50%  '
$c.Style = 'Normal'
$c = $ws.Range('D34')
$c.Value = "'" + '3.07'
$c.Style = 'Normal'
$c = $ws.Range('E34')
$c.Value = "'" + '  -2.30%  '
$c.Style = 'Normal'
$c = $ws.Range('D35')
$c.Value = "'" + '1.392.39'
$c.Style = 'Normal'
$c = $ws.Range('E35')
$c.Value = "'" + '  -0.62%  '
$c.Style = 'Normal'
$c = $ws.Range('E36')
$c.Value = "'" + '  +0.41%  '
$c.Style = 'Normal'
$c = $ws.Range('E37')
$c.Value = "'" + '  -3.33%  '
$c.Style = 'Normal'
$c = $ws.Range('D38')
$c.Value = "'" + '2.36'
$c.Style = 'Normal'
$c = $ws.Range('E38')
$c.Value = "'" + '  +0.28%  '
$c.Style = 'Normal'
$c = $ws.Range('D39')
$c.Value = "'" + '2.61'
$c.Style = 'Normal'
$c = $ws.Range('E39')
$c.Value = "'" + '  +4.92%  '
$c.Style = 'Normal'
$c = $ws.Range('D40')
$c.Value = "'" + '0.0166'
$c.Style = 'Normal'
$c = $ws.Range('E40')
$c.Value = "'" + '  -1.00%  '
$c.Style = 'Normal'
$c = $ws.Range('D41')
$c.Value = "'" + '0.530'
$c.Style = 'Normal'
$c = $ws.Range('E41')
$c.Value = "'" + '  -1.57%  '
$c.Style = 'Normal'
$c = $ws.Range('E42')
$c.Value = "'" + '  +0.05%  '
$c.Style = 'Normal'
$c = $ws.Range('E43')
$c.Value = "'" + '  +1.70%  '
$c.Style = 'Normal'
$c = $ws.Range('D44')
$c.Value = "'" + '0.787'
$c.Style = 'Normal'
$c = $ws.Range('E44')
$c.Value = "'" + '  -3.56%  '
$c.Style = 'Normal'
$c = $ws.Range('E45')
$c.Value = "'" + '  -1.53%  '
$c.Style = 'Normal'
$c = $ws.Range('E46')
$c.Value = "'" + '  -4.32%  '
$c.Style = 'Normal'
$c = $ws.Range('D47')
$c.Value = "'" + '62.59'
$c.Style = 'Normal'
$c = $ws.Range('E47')
$c.Value = "'" + '  -3.23%  '
$c.Style = 'Normal'
$c = $ws.Range('D48')
$c.Value = "'" + '1.705.63'
$c.Style = 'Normal'
$c = $ws.Range('E48')
$c.Value = "'" + '  -1.67%  '
$c.Style = 'Normal'
$c = $ws.Range('D49')
$c.Value = "'" + '85.94'
$c.Style = 'Normal'
$c = $ws.Range('E49')
$c.Value = "'" + '  -1.57%  '
$c.Style = 'Normal'
$c = $ws.Range('D50')
$c.Value = "'" + '0.0₆0101'
$c.Style = 'Normal'
$c = $ws.Range('E50')
$c.Value = "'" + '  -5.41%  '
$c.Style = 'Normal'
$c = $ws.Range('D51')
$c.Value = "'" + '0.0518'
$c.Style = 'Normal'
$c = $ws.Range('E51')
$c.Value = "'" + '  -2.12%  '
$c.Style = 'Normal'
